# "add score to ply func result"
# - add a new worksheet "add scoreply func 3.20" holding the (new) score-to-ply
#   results that previously lived in column H of "startimax 3.19"
# - clear column H on "startimax 3.19" (values removed, formatting kept where it
#   existed) since that data now lives on its own sheet
# - convert the ms timing row on "startimax 3.19" into seconds

$wb = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item(2)

# ---------------------------------------------------------------------------
# 1. Add the new sheet at the end of the workbook.
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws3 = $wb.Worksheets.Add($null, $lastSheet)
$ws3.Name = "add scoreply func 3.20"

# ---------------------------------------------------------------------------
# 2. Fill in the new sheet's data (mirrors the old column H layout, but with
#    its own result values).
# ---------------------------------------------------------------------------
$ws3.Range("B1").Value = "worst_rate:0.7 total:10 (add score to ply func) 30000:0 1000000:1 else 2 "

$ws3.Range("A2").Value = "128(%)"
$ws3.Range("B2").Value = 1

$ws3.Range("A3").Value = "256(%)"
$ws3.Range("B3").Value = 1

$ws3.Range("A4").Value = "512(%)"
$ws3.Range("B4").Value = 1

$ws3.Range("A5").Value = "1024(%)"
$ws3.Range("B5").Value = 1

$ws3.Range("A6").Value = "2048(%)"
$ws3.Range("B6").Value = 1

$ws3.Range("A7").Value = "4096(%)"
$ws3.Range("B7").Value = 0.8

$ws3.Range("A8").Value = "8192(%)"
$ws3.Range("B8").Value = 0.7

$ws3.Range("A9").Value = "16384(%)"
$ws3.Range("B9").Value = 0

$ws3.Range("A10").Value = "mean"
$ws3.Range("B10").Value = 123664

$ws3.Range("A11").Value = "max"
$ws3.Range("B11").Value = 178872

$ws3.Range("A12").Value = "time(sec)"
$ws3.Range("B12").Value = 2974.186

# Match the source percentages/number formatting used on the other sheets.
$ws3.Range("B2:B9").Style = "Percent"
$ws3.Range("B10").NumberFormat = $ws2.Range("F10").NumberFormat

$ws3.Columns.Item(1).ColumnWidth = 19.5
$ws3.Columns.Item(2).ColumnWidth = 55.5

# View: zoom + selection, then hand the active tab back to "startimax 3.19"
# (done at the very end, once everything else has been set up).
[void]$ws3.Activate()
$excel.ActiveWindow.Zoom = 108
$ws3.Range("B1").Select() | Out-Null

# ---------------------------------------------------------------------------
# 3. Column H on "startimax 3.19" now lives on the new sheet -- drop the
#    duplicated values. Rows that had formatting keep it (empty styled cell),
#    rows that had none (H1, H11, H12) disappear entirely.
# ---------------------------------------------------------------------------
$ws2.Range("H1").ClearContents()
$ws2.Range("H2:H10").ClearContents()
$ws2.Range("H11").ClearContents()
$ws2.Range("H12").ClearContents()

# ---------------------------------------------------------------------------
# 4. The "time(ms)" row becomes "time(sec)" with values divided by 1000.
# ---------------------------------------------------------------------------
$ws2.Range("A12").Value = "time(sec)"
$ws2.Range("D12").Value = 399.79199999999997
$ws2.Range("E12").Value = 6.8410000000000002
$ws2.Range("F12").Value = 240.01499999999999

# ---------------------------------------------------------------------------
# 5. Widen column H (it now only holds the header placeholder width) and
#    update the view: no more frozen/scrolled topLeftCell, new zoom + new
#    selection, and make "startimax 3.19" the active sheet/tab again.
# ---------------------------------------------------------------------------
$ws2.Columns.Item(8).ColumnWidth = 64.75

[void]$ws2.Activate()
$excel.ActiveWindow.Zoom = 106
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$ws2.Range("F12").Select() | Out-Null
